# "Update Venn diagram. Without borders."
#
# 1) Refresh the cached "datetimeFigureOut" field text on the slide master
#    and every slide layout (2018-03-22 -> 2018-04-14).
# 2) Give the (previously empty) 9th slide a fresh slide id — in the
#    original deck this happens because the slide content was rebuilt
#    (PowerPoint hands out the next free id, 267, replacing 266).
# 3) Draw the Venn diagram: three borderless ovals plus three text labels
#    (Monetizing / Engineering / Strategizing).

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($shapes, $text) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $text
            break
        }
    }
}

# --- 1) datetimeFigureOut field text -------------------------------------
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "2018-04-14"

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DatePlaceholderText $layouts.Item($i).Shapes "2018-04-14"
}

# --- 2) Renumber slide 9's sldId (266 -> 267) -----------------------------
$oldSlide = $p.Slides.Item(9)
$oldSlide.Copy()
$newRange = $p.Slides.Paste(9)
$newSlide = $newRange.Item(1)
# the original (still id 266) now sits one slot later; drop it
$p.Slides.Item(10).Delete()

# --- 3) Venn diagram shapes on the (new) slide 9 --------------------------
$ovalType = 9   # msoShapeOval
$rectType = 1   # msoShapeRectangle

$oval1 = $newSlide.Shapes.AddShape($ovalType, 230.19, 56.51622047244094, 283.3548031496063, 270.0)
$oval1.Fill.Visible = 0
$oval1.Line.Weight = 3
$oval1.Line.Visible = 0

$oval2 = $newSlide.Shapes.AddShape($ovalType, 338.7096850393701, 243.483937007874, 283.3548031496063, 270.0)
$oval2.Fill.Visible = 0
$oval2.Line.Weight = 3
$oval2.Line.Visible = 0

$oval3 = $newSlide.Shapes.AddShape($ovalType, 445.4156692913386, 56.51622047244094, 283.3548031496063, 270.0)
$oval3.Fill.Visible = 0
$oval3.Line.Weight = 3
$oval3.Line.Visible = 0

function Add-VennLabel($slide, $left, $top, $width, $height, $text, $languageId) {
    $shp = $slide.Shapes.AddShape($rectType, $left, $top, $width, $height)
    $tr = $shp.TextFrame.TextRange
    $tr.Text = $text
    $tr.Font.Size = 28
    $tr.Font.Name = "Arial"
    $tr.Font.NameComplexScript = "Arial"
    $tr.ParagraphFormat.Alignment = 2
    $tr.LanguageID = $languageId
    $shp.TextFrame.WordWrap = 0
    $shp.TextFrame.AutoSize = 1
    return $shp
}

Add-VennLabel $newSlide 544.4425984251968 150.31779527559056 151.6163779527559 41.198425196850394 "Monetizing" "sv-SE" | Out-Null
Add-VennLabel $newSlide 253.95629921259842 150.31779527559056 165.87929133858268 41.198425196850394 "Engineering" "sv-SE" | Out-Null
Add-VennLabel $newSlide 400.26637795275593 371.3848031496063 164.11220472440945 41.198425196850394 "Strategizing" "sv-SE" | Out-Null
